# TC10_Canine_Filter_StageOfDisease-5b.xlsx - "10 icdc scripts for jenkins"
#
# The FilesTab query (cell B4 on the "startup" sheet) is replaced with a
# simplified version of the Neo4j/Cypher query that drops the `File Type`
# and `Breed` output columns. The CasesTab (B2) and SamplesTab (B3) query
# text is unchanged. The sheet selection is moved from C2 to B4 to match
# where the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# New, shorter FilesTab Cypher query (File Type / Breed columns removed).
$newFilesQuery = @'

MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE diag.stage_of_disease IN ['Vb']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newFilesQuery

# Move the sheet's selection to the cell that was just edited.
$ws.Activate()
$ws.Range("B4").Select()

Write-Output "FilesTab query updated and selection moved to B4"
